$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (only B2 changes)
$ws.Range("B2").Value = 9891198044393078

# Row 3 - RandomForestRegressor (B3, C3, D3 change)
$ws.Range("B3").Value = 2593537460900.024
$ws.Range("C3").Value = 2499133910589.824
$ws.Range("D3").Value = 438183309651149.9

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor (A4, B4, C4, D4 change)
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2508509471790.415
$ws.Range("C4").Value = 2477383749676.61
$ws.Range("D4").Value = 151016259939687.9

# Row 5 - AdaBoostRegressor -> MLPRegressor (A5, B5, C5, D5 change)
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 312123613310573.9
$ws.Range("C5").Value = 396139552004133.2
$ws.Range("D5").Value = 3588734387100720
